# Registro_Proyectos_Actividades_ESTRATEGIAS.xlsx
# Sheet "Actividades": the ATEM_EST_2025 block grows from 4 activity rows
# (rows 2-5) to 6 activity rows (rows 2-7). Two new rows are inserted after
# the current row 5, pushing everything below (RECUP_APREND_EST_2024, ...,
# DOVE_EST_2025 blocks) down by two rows. The existing rows 2-5 also get
# updated content, and the two freshly inserted rows (6-7) are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades")

# Insert two new blank rows right before the current row 6 (shifts rows
# 6:33 down to 8:35, exactly matching the diff's dimension A1:K35).
$ws.Rows.Item(6).Resize(2).Insert()

# --- Row 2 (N=1, Encuentro de Formación Ciclo 1 año 1) ---------------------
# Only D2 (Número de Beneficiarios) changes.
$ws.Range("D2").Value = "218 Docentes `n"

# --- Row 3 (N=2, Visitas de acompañamiento situado) -------------------------
$ws.Range("B3").Value = "Visitas de acompañamiento situado ciclo 1 "
$ws.Range("D3").Value = "26 equipos gestores de la media  "
$ws.Range("E3").Value = "NO"
$ws.Range("F3").Value = ""

# --- Row 4 (N=3, previously blank -> Encuentro de Formación Ciclo 2) -------
$ws.Range("B4").Value = "Encuentro de Formación Ciclo 2 año 1 "
$ws.Range("C4").Value = "Docentes y Diectivos Docentes"
$ws.Range("D4").Value = "87 Docentes "
$ws.Range("E4").Value = "NO"

# --- Row 5 (N=4, previously blank -> Visitas de acompañamiento ciclo 2) ----
$ws.Range("B5").Value = "Visitas de acompañamiento situado ciclo 2"
$ws.Range("C5").Value = "Establecimientos Educativos"
$ws.Range("D5").Value = "26 equipos gestores de la media"
$ws.Range("E5").Value = "NO"

# --- Row 6 (new, N=5, Encuentro de Formación Ciclo 3 año 1) ----------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Encuentro de Formación Ciclo 3 año 1 "
$ws.Range("C6").Value = "Docentes y Diectivos Docentes"
$ws.Range("D6").Value = "134 Docentes "
$ws.Range("E6").Value = "NO"
$ws.Range("I6").Value = "Sin observaciones"
$ws.Range("J6").Value = "ATEM_EST_2025"
$ws.Range("K6").Value = "ALIANZA PARA LA TRANSFORMACIÓN DELA EDUCACIÓN MEDIA - ATEM"

# --- Row 7 (new, N=6, Visitas de acompañamiento situado ciclo 3) -----------
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Visitas de acompañamiento situado ciclo 3"
$ws.Range("C7").Value = "Establecimientos Educativos"
$ws.Range("D7").Value = "26 equipos gestores de la media"
$ws.Range("E7").Value = "NO"
$ws.Range("I7").Value = "Sin observaciones"
$ws.Range("J7").Value = "ATEM_EST_2025"
$ws.Range("K7").Value = "ALIANZA PARA LA TRANSFORMACIÓN DELA EDUCACIÓN MEDIA - ATEM"
